# ASCOM 0.3.0.0 - Updates
# - Added Slew rate support to MoveAxis
#
# This adds a new "Implemented in" (version) column E to the
# "ASCOM ItelescopeV3 Methods" sheet, tagging every already-implemented
# method with version 0.2.0.0, flips MoveAxis (row 12) from "No" to
# "Yes" now that it supports slew rate, and tags it with the new
# version 0.3.0.0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASCOM ItelescopeV3 Methods")

# MoveAxis now implemented, added in this release -> write this cell
# first so "0.3.0.0" lands before "0.2.0.0" in the shared-strings table.
$ws.Range("B12").Value = "Yes"
$ws.Range("E12").Value = "0.3.0.0"

# All the other already-implemented methods get tagged with the
# version they were originally implemented in (0.2.0.0).
$implementedRows = @(2, 3, 4, 5, 6, 7, 8, 10, 13, 14, 16, 19, 20, 21, 22, 24, 25, 26)
foreach ($r in $implementedRows) {
    $ws.Range("E$r").Value = "0.2.0.0"
}

# Match the cursor position recorded in the saved file.
$ws.Range("D5").Select()
